# Updated cryptos list on Thu May 23 13:17:01 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.427.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.58%  '

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.848.66'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.41%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.05%  '

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.65%  '

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.35%  '

# Row 7 - LidoStakedEther
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.847.85'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.46%  '

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.62%  '

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.163'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.02%  '

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.36'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.82%  '

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.477'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.23%  '

# Row 13 - Avalanche
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.45'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.25%  '

# Row 14 - ShibaInu
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000250'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.36%  '

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.506.37'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.84%  '

# Row 16 - WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.854.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.86%  '

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.433.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.48%  '

# Row 18 - Polkadot
$ws.Range("E18").Value = '  -0.76%  '

# Row 19 - TRON
$ws.Range("E19").Value = '  -3.44%  '

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.30%  '

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '500.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.34%  '

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.89%  '

# Row 23 - Polygon
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.744'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.15%  '

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.53'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.93%  '

# Row 25 - Fetch.AI
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.72%  '

# Row 26 - PEPE
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000138'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.95%  '

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.57%  '

# Row 28 - RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.40%  '

# Row 29 - Dai
$ws.Range("E29").Value = '  -0.15%  '

# Row 30 - ImmutableX
$ws.Range("E30").Value = '  +4.80%  '

# Row 31 - PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.97'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.06%  '

# Row 32 - EthereumClassic
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.89'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.46%  '

# Row 33 - NEARProtocol
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.93'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.46%  '

# Row 34 - Hedera
$ws.Range("E34").Value = '  +0.06%  '

# Row 35 - FirstDigitalUSD
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.03%  '

# Row 36 - Filecoin row becomes Mantle
$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.03'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.21%  '

# Row 37 - Mantle row becomes Filecoin
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.05%  '

# Row 38 - Kaspa
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.140'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.55%  '

# Row 39 - Bittensor
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '468.58'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.36%  '

# Row 40 - TheGraph
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.331'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.28%  '

# Row 41 - Stacks
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.05'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.06%  '

# Row 42 - OKB
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '49.49'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.93%  '

# Row 43 - dogwifhat
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.89'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.78%  '

# Row 44 - Cosmos
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.49'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.05%  '

# Row 45 - Arweave
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.76%  '

# Row 46 - Maker
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.895.22'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.49%  '

# Row 47 - VeChain
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0360'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.16%  '

# Row 48 - InjectiveProtocol row becomes Monero
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '139.45'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.63%  '

# Row 49 - Monero row becomes InjectiveProtocol
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '27.28'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.85%  '

# Row 50 - USDe
$ws.Range("E50").Value = '  +0.01%  '

# Row 51 - EnergySwap row becomes ThetaToken
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.38'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.15%  '
